$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The dataset used to live in A1:F14 (header in row 1, columns A-F).
# It now lives one column over and one row down, in B2:G15 - insert a
# blank column before A and a blank row before 1 to shift everything.
$ws.Columns("A:A").Insert()
$ws.Rows("1:1").Insert()

# Add a thin box border around every cell of the (now shifted) table,
# header row included.
$ws.Range("B2:G15").Borders.LineStyle = 1

# Page was set to print in portrait orientation.
$ws.PageSetup.Orientation = 1

# Selection left on F19 after the edits.
$null = $ws.Range("F19").Select()
